# Bot5 GUI beta terminada
# Clears the leftover "feriados" date/count values on the parametrosInicio
# sheet (B5: date "31.01.2023", B6: numeric 10) while leaving their cell
# styles intact, and updates the active selection to I12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parametrosInicio")

# Clear the values but keep formatting (style indexes s="2" / s="4" stay put).
$ws.Range("B5").ClearContents()
$ws.Range("B6").ClearContents()

# Move/store the sheet's active cell selection, as last left by the user.
$ws.Activate()
$ws.Range("I12").Select()
